$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new column before column A. Excel shifts the existing
#    PRODUCT_TYPE_ID / PRODUCT_MERK / PRODUCT_TYPE / CATEGORY columns (and
#    their formatting) one slot to the right, which is exactly what the
#    target layout needs (B:E instead of A:D).
# ---------------------------------------------------------------------------
$ws.Columns("A").Insert()

# ---------------------------------------------------------------------------
# 2. New header cell: "NOMOR"
# ---------------------------------------------------------------------------
$ws.Range("A1").Value2 = "NOMOR"

# ---------------------------------------------------------------------------
# 3. Format the new data cells A2:A8 with the same thin-box border already
#    used throughout the sheet (reuses the existing border style).
# ---------------------------------------------------------------------------
$dataRng = $ws.Range("A2:A8")
$dataRng.Borders.Color = 0
$dataRng.Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# 4. Format the new header cell A1: bold, centered, orange fill, no border.
# ---------------------------------------------------------------------------
$hdrA1 = $ws.Range("A1")
$hdrA1.Font.Bold = $true
$hdrA1.Interior.Color = 49407
$hdrA1.HorizontalAlignment = -4108
$hdrA1.VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 5. The rest of the header row (now B1:E1) loses the italic styling it used
#    to have, keeping bold -- i.e. bold-italic becomes plain bold.
# ---------------------------------------------------------------------------
$ws.Range("B1:E1").Font.Italic = $false

# ---------------------------------------------------------------------------
# 6. Column widths -- match the author's final widths as closely as this
#    engine's character-width rounding allows.
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 6.833333333333333
$ws.Columns("B").ColumnWidth = 19.666666666666668
$ws.Columns("C").ColumnWidth = 16.333333333333332
$ws.Columns("D").ColumnWidth = 15.666666666666666
$ws.Columns("E").ColumnWidth = 11.5

# ---------------------------------------------------------------------------
# 7. Selection, matching the saved sheet view in the target file.
# ---------------------------------------------------------------------------
$ws.Range("B10").Select()

Write-Output "Layout_ProductType: NOMOR column inserted"
